# Weekly refresh of the "Betarraga" price series: a new week's worth of
# data (Primera/Segunda) is inserted at the top of the data block (rows
# 220-221), pushing every existing week down by two rows. The two oldest
# rows that fall off the bottom of the original range are re-appended at
# the end (rows 262-263) so no historical data is lost.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 220; rows 220-261 shift down to 222-263.
$ws.Range("A220:A221").EntireRow.Insert()

# New row 220 - Betarraga, "Primera" quality, week of 2022-03-24 (serial 44644)
$ws.Range("A220").Value = 8
$ws.Range("B220").Value = "Terminal La Palmera de La Serena"
$ws.Range("C220").Value = "Coquimbo"
$ws.Range("D220").Value = 44644
$ws.Range("E220").Value = 4
$ws.Range("F220").Value = 100114014
$ws.Range("G220").Value = "Betarraga"
$ws.Range("H220").Value = "Sin especificar"
$ws.Range("I220").Value = "Primera"
$ws.Range("J220").Value = 2200
$ws.Range("K220").Value = 500
$ws.Range("L220").Value = 600
$ws.Range("M220").Value = 550
$ws.Range("N220").Value = "`$/paquete 3 unidades"
$ws.Range("O220").Value = "Provincia del Elquí"
$ws.Range("P220").Value = 183
$ws.Range("Q220").Value = 3
$ws.Range("R220").Value = "Hortaliza"

# New row 221 - Betarraga, "Segunda" quality, week of 2022-03-24 (serial 44644)
$ws.Range("A221").Value = 8
$ws.Range("B221").Value = "Terminal La Palmera de La Serena"
$ws.Range("C221").Value = "Coquimbo"
$ws.Range("D221").Value = 44644
$ws.Range("E221").Value = 4
$ws.Range("F221").Value = 100114014
$ws.Range("G221").Value = "Betarraga"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Segunda"
$ws.Range("J221").Value = 1460
$ws.Range("K221").Value = 400
$ws.Range("L221").Value = 450
$ws.Range("M221").Value = 425
$ws.Range("N221").Value = "`$/paquete 3 unidades"
$ws.Range("O221").Value = "Provincia del Elquí"
$ws.Range("P221").Value = 142
$ws.Range("Q221").Value = 3
$ws.Range("R221").Value = "Hortaliza"
